$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new volunteer time log rows (12 and 13) above the totals row.
$ws.Range("A12").Value = "11:36AM 12-28-2017"
$ws.Range("B12").Value = "12:38PM 12-28-2017"
$ws.Range("C12").Value = 62

$ws.Range("A13").Value = "1:54PM 12-12-2017"
$ws.Range("B13").Value = "5:42PM 12-28-2017"
$ws.Range("C13").Value = 228

# Widen column A slightly to fit the new data.
$ws.Columns.Item(1).ColumnWidth = 19

# Update the view: zoom in to 140% and move the selection to B14.
$excel.ActiveWindow.Zoom = 140
$ws.Range("B14").Select() | Out-Null
